# Update "想去人数" (interested-count) figures that changed between data pulls.
# Sheet "展览" (sheetId 1 / sheet1.xml)
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 3583
$wsExhibit.Range("F5").Value = 3583
$wsExhibit.Range("F7").Value = 5113
$wsExhibit.Range("F8").Value = 5113
$wsExhibit.Range("F9").Value = 532
$wsExhibit.Range("F12").Value = 693
$wsExhibit.Range("F14").Value = 94
$wsExhibit.Range("F18").Value = 37
$wsExhibit.Range("F23").Value = 4916
$wsExhibit.Range("F24").Value = 4916
$wsExhibit.Range("F38").Value = 1025
$wsExhibit.Range("F42").Value = 874
$wsExhibit.Range("F43").Value = 999
$wsExhibit.Range("F44").Value = 2030

# Sheet "全部类型" (sheetId 4 / sheet4.xml) mirrors the same events at different rows.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 3583
$wsAll.Range("F8").Value = 3583
$wsAll.Range("F10").Value = 5113
$wsAll.Range("F11").Value = 5113
$wsAll.Range("F12").Value = 532
$wsAll.Range("F15").Value = 693
$wsAll.Range("F17").Value = 94
$wsAll.Range("F21").Value = 37
$wsAll.Range("F27").Value = 4916
$wsAll.Range("F28").Value = 4916
$wsAll.Range("F43").Value = 1025
$wsAll.Range("F47").Value = 874
$wsAll.Range("F48").Value = 999
$wsAll.Range("F50").Value = 2030
